$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Near the end of the document: drop the duplicated bold title
#    paragraph and rewrite the italic "meta description" paragraph
#    into the new image-prompt text. Do this first (before the
#    insertion below) so paragraph indices stay simple/predictable.
# ---------------------------------------------------------------------
$boldTitlePara = $d.Paragraphs.Item(55)
$descParaBefore = $d.Paragraphs.Item(56)

# sanity checks (raise if the document doesn't look as expected)
if ($boldTitlePara.Range.Text.TrimEnd() -ne "Play Chicken Fox Free: Cute Farm-Inspired Slot Game") {
    throw "Unexpected paragraph 55 content: $($boldTitlePara.Range.Text)"
}
if ($descParaBefore.Range.Text.TrimEnd() -ne "Check out our review of Chicken Fox, a cute and playful online slot game with exciting bonus features. Play for free and enjoy the farm-inspired graphics.") {
    throw "Unexpected paragraph 56 content: $($descParaBefore.Range.Text)"
}

# Delete the whole duplicated bold-title paragraph, including its
# trailing paragraph mark.
$delRange = $d.Range($boldTitlePara.Range.Start, $descParaBefore.Range.Start)
$delRange.Delete()

# The italic paragraph is now paragraph 55; rewrite its text in place
# (keeps the run's rPr -> <w:i/> and the paragraph's leading empty run).
$descPara = $d.Paragraphs.Item(55)
$descRange = $d.Range($descPara.Range.Start, $descPara.Range.End)
$descRange.Text = "Create a feature image for Chicken Fox that captures the cute farmyard theme of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be depicted alongside the game's farm animals, such as chickens, goats, and pigs, with a background of a traditional farm setting, complete with a barn and rolling hills. The image should convey the fun, playful nature of the game, while also highlighting its unique features, such as the Free Games feature and multipliers."

# ---------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Build the paragraph's runs (leading empty run + bold "Meta description"
# run) via raw OOXML so the shape matches the rest of the document
# (every body paragraph here starts with an empty <w:r/>).
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

$metaPara2 = $d.Paragraphs.Item(2)
$metaTail = $metaPara2.Range
$metaTail.InsertAfter(": Check out our review of Chicken Fox, a cute and playful online slot game with exciting bonus features. Play for free and enjoy the farm-inspired graphics.")
